$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G5").Value = "93-100%"
$ws.Range("H5").Value = "56-86%"
$ws.Range("I5").Value = "93-100%"

$ws.Range("F6").Value = "83-100%"
$ws.Range("G6").Value = "86-93%"
$ws.Range("H6").Value = "30-70%"
$ws.Range("I6").Value = 1
$ws.Range("I6").NumberFormat = "0%"

$ws.Range("F7").Value = "67-90%"
$ws.Range("G7").Value = "70-90%"
$ws.Range("H7").Value = "77-86%"
$ws.Range("I7").Value = "86-100%"

$ws.Range("F8").Value = "36-70%"
$ws.Range("G8").Value = "50-80%"
$ws.Range("H8").Value = "60-97%"
$ws.Range("I8").Value = "70-100%"
$ws.Range("I8").NumberFormat = "0%"

$ws.Range("K18").Select()
